# GFG-Deletion from a Circular Linked List
# Append a new row (row 19) to the Linked List question tracker sheet:
#   A19 = "GFG"   (reuses the existing "GFG" shared string)
#   B19 = "Deletion from a Circular Linked List"  (new shared string)
# and move the selection to the newly-added cell, matching the author's
# final cursor position after data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "GFG"
$ws.Range("B19").Value = "Deletion from a Circular Linked List"

# Leave the selection on the new cell, same as the source workbook.
$ws.Range("B19").Select()
